$d = $word.ActiveDocument

# --- Edit 1: "Programa resumido" (short/summary) English italic paragraph ---
$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$findText1 = "1) Conceptual basis for transport phenomena study2) General properties of fluids 3) Kinematics of fluids.4) Conservation Equations in Integral form.5) Differential Equations of Fluid Flow. 6) Boundary Layer Theory.7) Flow in ducts:"
$replaceText1 = "1) Conceptual basis for transport phenomena study^l2) General properties of fluids ^l3) Kinematics of fluids.^l4) Conservation Equations in Integral form.^l5) Differential Equations of Fluid Flow. ^l6) Boundary Layer Theory.^l7) Flow in ducts:"
$found1 = $rng1.Find.Execute($findText1, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText1, 2)
Write-Host "Edit 1 found/replaced: $found1"
if (-not $found1) {
    throw "Edit 1: target text for 'Programa resumido' English paragraph was not found."
}

# --- Edit 2: "Programa" (detailed) English italic paragraph ---
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$findText2 = "1) Conceptual basis for transport phenomena studyFluids and the continuous hypothesis. Importance of dimensional analysis and use of dimensionless numbers. Basic laws for mass, heat and motion amount transfer. General law for transport phenomena. Molecular diffusivity, thermal conductivity and viscosity. Simultaneous transport of mass, heat and motion amount. Integral and differential formulation.2) General properties of fluids: Specific mass, specific weight, specific volume. Tension and pressure. Newtonian and non-Newtonian fluids. Viscosity. Surface tension and capillarity. Volumetric elasticity modulus and compressibility.3) Fluid Kinematics: Description of a Fluid Motion: Euler and Lagrange method - Field of fluid flow- Permanent and transient flow - Trajectories and Streamlines - System and volume control – Unidimensional and bidimensional flows. Uniform flow. Laminar and turbulent flow: Reynolds number.4) Conservation Equations in Integral form: Flow of a magnitude. Mass conservation, continuity. Specific forms for the integral expression. Amount conservation of linear motion. Energy conservation. Bernoulli Equation. Applications.5) Differential Equations of Fluid Flow: Mass conservation Equation and continuity. Energy equations. Navier-Stokes equations. Applications."
$replaceText2 = "1) Conceptual basis for transport phenomena study^lFluids and the continuous hypothesis. Importance of dimensional analysis and use of dimensionless numbers. Basic laws for mass, heat and motion amount transfer. General law for transport phenomena. Molecular diffusivity, thermal conductivity and viscosity. Simultaneous transport of mass, heat and motion amount. Integral and differential formulation.^l2) General properties of fluids: Specific mass, specific weight, specific volume. Tension and pressure. Newtonian and non-Newtonian fluids. Viscosity. Surface tension and capillarity. Volumetric elasticity modulus and compressibility.^l3) Fluid Kinematics: Description of a Fluid Motion: Euler and Lagrange method - Field of fluid flow- Permanent and transient flow - Trajectories and Streamlines - System and volume control – Unidimensional and bidimensional flows. Uniform flow. Laminar and turbulent flow: Reynolds number.^l4) Conservation Equations in Integral form: Flow of a magnitude. Mass conservation, continuity. Specific forms for the integral expression. Amount conservation of linear motion. Energy conservation. Bernoulli Equation. Applications.^l5) Differential Equations of Fluid Flow: Mass conservation Equation and continuity. Energy equations. Navier-Stokes equations. Applications."
$found2 = $rng2.Find.Execute($findText2, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText2, 2)
Write-Host "Edit 2 found/replaced: $found2"
if (-not $found2) {
    throw "Edit 2: target text for 'Programa' English paragraph was not found."
}
